$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Stats" column (H) values for the three NFT rows so that the
# property name is a quoted JSON-style string ("Strenght") instead of a
# bare word (Strenght), matching the corrected data export format.
$ws.Range("H2").Value = ' [["Strenght", 10, 100], ["Age", 1, 99]]'
$ws.Range("H3").Value = ' [["Strenght", 10, 100]]'
$ws.Range("H4").Value = ' ["Strenght", 10, 100]'

# Update the active selection on the sheet to H7, matching the saved
# selection state in the workbook.
$ws.Range("H7").Select()
